# Update 想去人数 (column F) values across sheets per source data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 2265
$ws.Cells.Item(3, 6).Value = 320
$ws.Cells.Item(4, 6).Value = 175
$ws.Cells.Item(5, 6).Value = 178
$ws.Cells.Item(6, 6).Value = 328
$ws.Cells.Item(8, 6).Value = 686
$ws.Cells.Item(9, 6).Value = 505
$ws.Cells.Item(10, 6).Value = 641
$ws.Cells.Item(12, 6).Value = 63
$ws.Cells.Item(13, 6).Value = 351
$ws.Cells.Item(14, 6).Value = 963
$ws.Cells.Item(15, 6).Value = 5
$ws.Cells.Item(16, 6).Value = 240
$ws.Cells.Item(17, 6).Value = 134
$ws.Cells.Item(20, 6).Value = 48
$ws.Cells.Item(22, 6).Value = 243
$ws.Cells.Item(23, 6).Value = 96

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 80
$ws.Cells.Item(6, 6).Value = 174
$ws.Cells.Item(7, 6).Value = 205
$ws.Cells.Item(8, 6).Value = 2564
$ws.Cells.Item(13, 6).Value = 25
$ws.Cells.Item(16, 6).Value = 2425

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 252
$ws.Cells.Item(3, 6).Value = 21
$ws.Cells.Item(4, 6).Value = 355

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 252
$ws.Cells.Item(3, 6).Value = 21
$ws.Cells.Item(6, 6).Value = 2265
$ws.Cells.Item(7, 6).Value = 355
$ws.Cells.Item(8, 6).Value = 320
$ws.Cells.Item(9, 6).Value = 175
$ws.Cells.Item(10, 6).Value = 178
$ws.Cells.Item(11, 6).Value = 328
$ws.Cells.Item(14, 6).Value = 80
$ws.Cells.Item(15, 6).Value = 174
$ws.Cells.Item(17, 6).Value = 686
$ws.Cells.Item(18, 6).Value = 505
$ws.Cells.Item(19, 6).Value = 641
$ws.Cells.Item(21, 6).Value = 63
$ws.Cells.Item(22, 6).Value = 351
$ws.Cells.Item(23, 6).Value = 963
$ws.Cells.Item(24, 6).Value = 205
$ws.Cells.Item(25, 6).Value = 2564
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(31, 6).Value = 25
$ws.Cells.Item(32, 6).Value = 240
$ws.Cells.Item(33, 6).Value = 134
$ws.Cells.Item(38, 6).Value = 48
$ws.Cells.Item(40, 6).Value = 243
$ws.Cells.Item(41, 6).Value = 96
$ws.Cells.Item(42, 6).Value = 2425
